$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "number_stimulus"
$ws.Range("B1").Value = "letter_stimulus"
$ws.Range("C1").Value = "cor_number_resp"
$ws.Range("D1").Value = "cor_par_resp"

# Data rows 2-9
$data = @(
    @(2, "A", "j", "f"),
    @(3, "B", "k", "g"),
    @(4, "C", "j", "g"),
    @(5, "E", "k", "f"),
    @(6, "I", "j", "f"),
    @(7, "L", "k", "g"),
    @(8, "U", "j", "f"),
    @(9, "W", "k", "g")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Data rows 10-17 (repeat of rows 2-9)
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Column widths: size columns A and B to fit their new (longer) content
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 14

# Selection
$ws.Range("I16").Select()
